$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.354.28"
$ws.Range("E2").Value = "  -2.80%  "
$ws.Range("D3").Value = "1.941.14"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7189"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -7.39%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3347"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.74"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07330"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8166"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08157"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.939.96"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.499"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.38"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").Value = "30.370.48"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008374"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +6.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.867"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "2.195.44"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.973"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.849"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.405"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1315"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -10.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.578"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.344"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.493"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.265"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05288"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.271"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7668"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02001"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.844"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.563"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4579"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.033"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8481"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.853"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.454"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.17"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4188"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.504"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.15%  "
